# Auto-generated edit script: updates cached numeric values in the
# "Leve profit" sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match
# a refreshed data pull from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H51").Value = 3364.7693
$ws.Range("J51").Value = 2142
$ws.Range("L51").Value = 2142
$ws.Range("N51").Value = -3110
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").Value = $null
$ws.Range("H62").Value = 4799.9443
$ws.Range("J62").Value = 5415.6665
$ws.Range("L62").Value = 5415.6665
$ws.Range("N62").Value = -6663.6665
$ws.Range("H65").Value = 4799.9443
$ws.Range("J65").Value = 5415.6665
$ws.Range("L65").Value = 27078.3325
$ws.Range("N65").Value = -33318.3325
$ws.Range("H68").Value = 124999
$ws.Range("J68").Value = 124999
$ws.Range("L68").Value = 124999
$ws.Range("N68").Value = -126497
$ws.Range("H71").Value = 124999
$ws.Range("J71").Value = 124999
$ws.Range("L71").Value = 374997
$ws.Range("N71").Value = -382485
$ws.Range("H113").Value = 2610.4
$ws.Range("I113").Value = 2593.4666
$ws.Range("K113").Value = 2593.4666
$ws.Range("M113").Value = 660.5333999999998
$ws.Range("H126").Value = 142890
$ws.Range("J126").Value = 142890
$ws.Range("L126").Value = 142890
$ws.Range("N126").Value = -152770
$ws.Range("H135").Value = 1078.2174
$ws.Range("I135").Value = 1042.421
$ws.Range("J135").Value = 1248.25
$ws.Range("K135").Value = 9381.789000000001
$ws.Range("L135").Value = 11234.25
$ws.Range("M135").Value = -6846.789000000001
$ws.Range("N135").Value = -16304.25
$ws.Range("H137").Value = 1881240.6
$ws.Range("J137").Value = 3270827.2
$ws.Range("L137").Value = 9812481.600000001
$ws.Range("N137").Value = -9817581.600000001
$ws.Range("H141").Value = 5221.1924
$ws.Range("I141").Value = 4739.6665
$ws.Range("K141").Value = 14218.9995
$ws.Range("M141").Value = -9038.999500000002

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H45").Value = 4849.227
$ws.Range("J45").Value = 4916.5
$ws.Range("L45").Value = 4916.5
$ws.Range("N45").Value = -5670.5
$ws.Range("H74").Value = 2264.0312
$ws.Range("I74").Value = 2118.75
$ws.Range("K74").Value = 2118.75
$ws.Range("M74").Value = -1244.75
$ws.Range("H77").Value = 2264.0312
$ws.Range("I77").Value = 2118.75
$ws.Range("K77").Value = 10593.75
$ws.Range("M77").Value = -6225.75
$ws.Range("H122").Value = 3946.682
$ws.Range("I122").Value = 3946.682
$ws.Range("K122").Value = 11840.046
$ws.Range("M122").Value = -9390.045999999998

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H124").Value = 75269
$ws.Range("J124").Value = 75269
$ws.Range("L124").Value = 75269
$ws.Range("N124").Value = -85089
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H22").Value = 10776.154
$ws.Range("I22").Value = 1099
$ws.Range("J22").Value = 32549.75
$ws.Range("K22").Value = 1099
$ws.Range("L22").Value = 32549.75
$ws.Range("M22").Value = -749
$ws.Range("N22").Value = -33249.75
$ws.Range("H31").Value = 3345.5962
$ws.Range("I31").Value = 2185.4614
$ws.Range("J31").Value = 3732.3076
$ws.Range("K31").Value = 2185.4614
$ws.Range("L31").Value = 3732.3076
$ws.Range("M31").Value = -1890.4614
$ws.Range("N31").Value = -4322.3076
$ws.Range("H34").Value = 3345.5962
$ws.Range("I34").Value = 2185.4614
$ws.Range("J34").Value = 3732.3076
$ws.Range("K34").Value = 2185.4614
$ws.Range("L34").Value = 3732.3076
$ws.Range("M34").Value = -1983.4614
$ws.Range("N34").Value = -4136.3076

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H68").Value = 1967.3043
$ws.Range("J68").Value = 2251.6
$ws.Range("L68").Value = 6754.799999999999
$ws.Range("N68").Value = -8376.799999999999
$ws.Range("H71").Value = 1967.3043
$ws.Range("J71").Value = 2251.6
$ws.Range("L71").Value = 20264.4
$ws.Range("N71").Value = -28376.4
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = $null
$ws.Range("N92").Value = $null
$ws.Range("H117").Value = 334410
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 334410
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 1003230
$ws.Range("M117").Value = $null
$ws.Range("N117").Value = -1010114
$ws.Range("H129").Value = 1210
$ws.Range("I129").Value = 613.3333
$ws.Range("K129").Value = 1839.9999
$ws.Range("M129").Value = 3160.0001

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 4499.9443
$ws.Range("J70").Value = 4499.9287
$ws.Range("L70").Value = 4499.9287
$ws.Range("N70").Value = -5039.9287
$ws.Range("H73").Value = 4499.9443
$ws.Range("J73").Value = 4499.9287
$ws.Range("L73").Value = 4499.9287
$ws.Range("N73").Value = -6371.9287
$ws.Range("H127").Value = 99990
$ws.Range("J127").Value = 99990
$ws.Range("L127").Value = 99990
$ws.Range("N127").Value = -109910
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null
$ws.Range("H140").Value = 110833.336
$ws.Range("J140").Value = 146250
$ws.Range("L140").Value = 146250
$ws.Range("N140").Value = -156610

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H46").Value = 7193.92
$ws.Range("J46").Value = 7929.591
$ws.Range("L46").Value = 7929.591
$ws.Range("N46").Value = -8305.591
$ws.Range("H68").Value = 1757.6666
$ws.Range("I68").Value = 1802.375
$ws.Range("J68").Value = 1400
$ws.Range("K68").Value = 1802.375
$ws.Range("L68").Value = 1400
$ws.Range("M68").Value = -1053.375
$ws.Range("N68").Value = -2898
$ws.Range("H71").Value = 1757.6666
$ws.Range("I71").Value = 1802.375
$ws.Range("J71").Value = 1400
$ws.Range("K71").Value = 9011.875
$ws.Range("L71").Value = 7000
$ws.Range("M71").Value = -5267.875
$ws.Range("N71").Value = -14488
$ws.Range("H119").Value = 95141
$ws.Range("J119").Value = 95141
$ws.Range("L119").Value = 95141
$ws.Range("N119").Value = -104817
$ws.Range("H132").Value = 628319.3
$ws.Range("I132").Value = 772317.25
$ws.Range("J132").Value = 4328.3335
$ws.Range("K132").Value = 2316951.75
$ws.Range("L132").Value = 12985.0005
$ws.Range("M132").Value = -2314421.75
$ws.Range("N132").Value = -18045.0005
$ws.Range("H141").Value = 599999
$ws.Range("J141").Value = 599999
$ws.Range("L141").Value = 599999
$ws.Range("N141").Value = -610359

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H62").Value = 7457.4287
$ws.Range("I62").Value = 8834
$ws.Range("J62").Value = 6425
$ws.Range("K62").Value = 8834
$ws.Range("L62").Value = 6425
$ws.Range("M62").Value = -8210
$ws.Range("N62").Value = -7673
$ws.Range("H65").Value = 7457.4287
$ws.Range("I65").Value = 8834
$ws.Range("J65").Value = 6425
$ws.Range("K65").Value = 44170
$ws.Range("L65").Value = 32125
$ws.Range("M65").Value = -41050
$ws.Range("N65").Value = -38365
$ws.Range("H113").Value = 724.6923
$ws.Range("I113").Value = 724.6923
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2174.0769
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -4.076900000000023
$ws.Range("N113").Value = $null
$ws.Range("H122").Value = 4911.6787
$ws.Range("I122").Value = 4961.08
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 14883.24
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -12433.24
$ws.Range("N122").Value = -18400
$ws.Range("H131").Value = 149984.5
$ws.Range("J131").Value = 149984.5
$ws.Range("L131").Value = 149984.5
$ws.Range("N131").Value = -160064.5
$ws.Range("H135").Value = 203476.67
$ws.Range("J135").Value = 203476.67
$ws.Range("L135").Value = 203476.67
$ws.Range("N135").Value = -213616.67
$ws.Range("H136").Value = 57053.473
$ws.Range("I136").Value = 4000.5
$ws.Range("J136").Value = 95637.45
$ws.Range("K136").Value = 12001.5
$ws.Range("L136").Value = 286912.35
$ws.Range("M136").Value = -9451.5
$ws.Range("N136").Value = -292012.35
$ws.Range("H141").Value = 176999.5
$ws.Range("J141").Value = 176999.5
$ws.Range("L141").Value = 176999.5
$ws.Range("N141").Value = -187359.5
